$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.792.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.707.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9970"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9965"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3934"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.499"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9959"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08826"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +11.07%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.134"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001362"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.701.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07171"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.304"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9966"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.39"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.799.61"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.021"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.337"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.72"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.924"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +15.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.604"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "145.40"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.887.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08830"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.176"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.070"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.231"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.03120"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2821"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8536"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +10.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.93"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09223"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.20"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.475"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.74"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.714"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7520"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.282"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.392"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9959"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08268"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.69%  "
